$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the last-changed date serial for every data
# row (rows 2-309). The update bumps that date by one day, from serial
# 46061 (2026-02-08) to 46062 (2026-02-09), for every row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$ws.Range("C2:C$lastRow").Value = 46062
